$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.175.00"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.784.30"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.02%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.36%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "316.29"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("E6").Value = "  +0.26%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5313"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -3.20%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3745"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -3.11%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07465"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.44%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.52"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.34%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.091"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.75%  "

$ws.Range("E12").Value = "  +0.25%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.38"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.71%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.090"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.88%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.222"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.72%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.766.67"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.70%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "88.91"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.22%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001055"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.47%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06479"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.21%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.43"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.94%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.913"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.12%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.226.21"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("E24").Value = "  -3.60%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.091"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.58%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "157.74"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.28%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.25"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.06%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.980.88"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.26%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.286"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -5.87%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "120.78"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.36%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.094"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.20%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.1042"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.61%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.662"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.18%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.509"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.02%  "

$ws.Range("E35").Value = "  -2.89%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.06367"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.13%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02272"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.08%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.984"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.80%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "8.463"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.98%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.6159"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.65%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "11.00"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -5.38%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.429"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.46%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.177"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("E44").Value = "  +0.21%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.29"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.672"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.29%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5757"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.86%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "125.79"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("E49").Value = "  +4.70%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.926"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.59%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06835"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.94%  "
